# Update the build/version string embedded throughout the workbook, for the
# new "Coal Mine Boundaries and Methane Sources - version 1.0.0" release.
#
# Old:  mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)
# New:  Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: <version string>"
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text containing the version string (single-quoted
# inside the sentence).
$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Grosvenor Coal Mine, Australia, M0045, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$aboutSheet.Range("A6").Value = $newCitation

# Column S ("build_version") on the data sheet holds the same version string
# for every data row (rows 2..last). Replace it wherever it matches the old
# value.
$lastRow = $dataSheet.Cells.Item($dataSheet.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
